$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values (participant/column identifiers)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) updated meanEMG legmaxROM values
$ws.Range("B2").Value = 235.16656560964404
$ws.Range("C2").Value = 172.67528433332257
$ws.Range("D2").Value = 235.73520860205826
$ws.Range("E2").Value = 174.58653062100834

# Row 3 (STR) updated meanEMG legmaxROM values
$ws.Range("B3").Value = 215.66093203200802
$ws.Range("C3").Value = 171.65790081672372
$ws.Range("D3").Value = 211.91721311210463
$ws.Range("E3").Value = 180.64958327106208

# Update selection to match new highlighted range
$ws.Range("B1:E3").Select()
